$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Receptor avg/total expr values, specificities, Edge weights/specificities)
$ws.Range("M2").Value = 0.000144
$ws.Range("N2").Value = 0.000432
$ws.Range("O2").Value = 0.6050420168067226
$ws.Range("P2").Value = 0.6050420168067226
$ws.Range("Q2").Value = 0.00007878523200000001
$ws.Range("R2").Value = 0.0007090670879999999
$ws.Range("S2").Value = 0.6050420168067226
$ws.Range("T2").Value = 0.6050420168067226

# Row 3 updates (Receptor-expressing cells, Receptor detection rate, and downstream values)
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.00009400000000000001
$ws.Range("N3").Value = 0.000282
$ws.Range("O3").Value = 0.3949579831932774
$ws.Range("P3").Value = 0.3949579831932774
$ws.Range("Q3").Value = 0.00005142924866666667
$ws.Range("R3").Value = 0.0004628632380000001
$ws.Range("S3").Value = 0.3949579831932774
$ws.Range("T3").Value = 0.3949579831932774
